## Adds a new "2021" data column (L) to the table, mirroring the existing
## D:K year columns, and moves the active selection as in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# L3: empty border cell, same format as the other cells in the divider
# row (D3:H3, J3) -- font 1 (Times New Roman 10), bottom medium border,
# vertical-center alignment.
# ---------------------------------------------------------------------
$ws.Range("D3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# L4: header "2021" -- bold font (fontId 2) with a bottom medium border,
# no alignment override. Start from K4's format (bold font + top+bottom
# border) and drop the top edge so the border collapses onto the
# existing "bottom-only" border definition instead of minting a new one.
# ---------------------------------------------------------------------
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Borders.Item(8).LineStyle = -4142
$ws.Range("L4").Value2 = 2021

# ---------------------------------------------------------------------
# L5: data value 2.3 -- regular (non-bold) font, no border.
# ---------------------------------------------------------------------
$ws.Range("K4").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Font.Bold = $false
$ws.Range("L5").Borders.Item(8).LineStyle = -4142
$ws.Range("L5").Borders.Item(9).LineStyle = -4142
$ws.Range("L5").Value2 = 2.2999999999999998

# ---------------------------------------------------------------------
# L6: data value 1.3 -- regular (non-bold) font with a bottom medium
# border (table's closing row).
# ---------------------------------------------------------------------
$ws.Range("K4").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Font.Bold = $false
$ws.Range("L6").Borders.Item(8).LineStyle = -4142
$ws.Range("L6").Value2 = 1.3

$excel.CutCopyMode = 0

# Match the cursor position left behind in the saved workbook.
$ws.Range("O5").Select() | Out-Null
